# Complete 0.2.0 changelog + site consistency
# Find/replace 0.1.6* -> 0.2.0 and related "Moved to ..." / bugfix link notes
# in the "Source table" sheet. Downstream sheets (DIMR/RR/FM mkdocs tables)
# consume these values through formulas referencing 'Source table', so they
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source table")

# --- Version bump: 0.1.6* -> 0.2.0 -------------------------------------
$versionCells = @("D3", "D22", "D40", "D41", "D43", "D44", "D62", "D63")
foreach ($cellRef in $versionCells) {
    $ws.Range($cellRef).Value = "0.2.0"
}

# --- "Moved to io.crosssections in 0.1.6." -> "...0.2.0" ---------------
$ws.Range("G31").Value = "Moved to io.crosssections in 0.2.0"

# --- "Moved to io.rr in 0.1.6" -> "...0.2.0" ----------------------------
$ws.Range("G59").Value = "Moved to io.rr in 0.2.0"
$ws.Range("G60").Value = "Moved to io.rr in 0.2.0"

# --- Critical bugfix note: fix markdown link syntax ---------------------
$ws.Range("G3").Value = "Critical bugfix for [#127](https://github.com/Deltares/HYDROLIB-core/issues/127)."

# --- Restore view state on the Source table sheet -----------------------
$ws.Select()
$ws.Range("G4").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
